# This script updates the cryptocurrency price/volume table with a fresh
# scrape of data (prices in column D, 1h volume % change in column E, and
# for row 51 the coin name/link in columns B/C as well).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All of these cells originally hold plain text (not real numbers), e.g.
# "207.38" or "1.557.28" are scraped price strings, and "  +0.86%  " is a
# padded percentage string. Assigning such strings via .Value/.Value2 makes
# Excel silently reinterpret anything that parses as a plain decimal (like
# "207.29") as a real number. Set-TextValue guards against that: for values
# that look like a plain number it temporarily forces a text format, writes
# the literal text, then restores the cell to the default "Normal" style so
# no stray number formatting is left behind. Values that are not plain
# numbers (e.g. "27.012.05", "  +0.93%  ", "BabyDogeCoin") are unaffected by
# this auto-conversion and are written directly.
function Set-TextValue($cell, $text) {
    if ($text -match '^-?[0-9]+(\.[0-9]+)?$') {
        $cell.NumberFormat = "@"
        $cell.Value2 = $text
        $cell.Style = "Normal"
    } else {
        $cell.Value2 = $text
    }
}

Set-TextValue $ws.Range("D2") "27.012.05"
Set-TextValue $ws.Range("E2") "  +0.93%  "
Set-TextValue $ws.Range("D3") "1.557.26"
Set-TextValue $ws.Range("E3") "  +0.68%  "
Set-TextValue $ws.Range("D5") "207.29"
Set-TextValue $ws.Range("E5") "  +0.68%  "
Set-TextValue $ws.Range("E6") "  +1.37%  "
Set-TextValue $ws.Range("E7") "  +0.32%  "
Set-TextValue $ws.Range("D8") "21.68"
Set-TextValue $ws.Range("E8") "  +1.47%  "
Set-TextValue $ws.Range("E9") "  +1.28%  "
Set-TextValue $ws.Range("E10") "  +1.59%  "
Set-TextValue $ws.Range("E11") "  +0.96%  "
Set-TextValue $ws.Range("D12") "1.779.31"
Set-TextValue $ws.Range("E12") "  +0.68%  "
Set-TextValue $ws.Range("D13") "1.557.46"
Set-TextValue $ws.Range("E13") "  +0.67%  "
Set-TextValue $ws.Range("E14") "  +1.67%  "
Set-TextValue $ws.Range("D15") "0.516"
Set-TextValue $ws.Range("E15") "  +1.10%  "
Set-TextValue $ws.Range("D16") "27.005.23"
Set-TextValue $ws.Range("E16") "  +0.91%  "
Set-TextValue $ws.Range("D17") "61.98"
Set-TextValue $ws.Range("E17") "  +1.46%  "
Set-TextValue $ws.Range("D18") "216.09"
Set-TextValue $ws.Range("E18") "  +1.32%  "
Set-TextValue $ws.Range("E19") "  +0.63%  "
Set-TextValue $ws.Range("E20") "  +0.59%  "
Set-TextValue $ws.Range("E21") "  +0.30%  "
Set-TextValue $ws.Range("D22") "4.03"
Set-TextValue $ws.Range("D23") "9.22"
Set-TextValue $ws.Range("E23") "  +2.90%  "
Set-TextValue $ws.Range("D24") "1.97"
Set-TextValue $ws.Range("E24") "  -0.85%  "
Set-TextValue $ws.Range("D25") "152.47"
Set-TextValue $ws.Range("E25") "  -0.46%  "
Set-TextValue $ws.Range("E27") "  +0.26%  "
Set-TextValue $ws.Range("E28") "  +0.33%  "
Set-TextValue $ws.Range("E29") "  +1.55%  "
Set-TextValue $ws.Range("E30") "  +0.58%  "
Set-TextValue $ws.Range("D31") "1.09"
Set-TextValue $ws.Range("E31") "  -0.73%  "
Set-TextValue $ws.Range("E32") "  +1.25%  "
Set-TextValue $ws.Range("D33") "1.401.46"
Set-TextValue $ws.Range("E33") "  +4.14%  "
Set-TextValue $ws.Range("E34") "  +3.34%  "
Set-TextValue $ws.Range("E35") "  +3.87%  "
Set-TextValue $ws.Range("D36") "0.965"
Set-TextValue $ws.Range("E36") "  +4.37%  "
Set-TextValue $ws.Range("E37") "  +0.15%  "
Set-TextValue $ws.Range("E38") "  +1.58%  "
Set-TextValue $ws.Range("D39") "0.523"
Set-TextValue $ws.Range("E39") "  +1.03%  "
Set-TextValue $ws.Range("E40") "  +1.47%  "
Set-TextValue $ws.Range("E41") "  +0.36%  "
Set-TextValue $ws.Range("D42") "0.990"
Set-TextValue $ws.Range("E42") "  -0.23%  "
Set-TextValue $ws.Range("D43") "2.27"
Set-TextValue $ws.Range("E43") "  +3.68%  "
Set-TextValue $ws.Range("E44") "  -3.50%  "
Set-TextValue $ws.Range("D45") "63.93"
Set-TextValue $ws.Range("E45") "  +1.79%  "
Set-TextValue $ws.Range("E46") "  -0.11%  "
Set-TextValue $ws.Range("D47") "1.692.77"
Set-TextValue $ws.Range("E47") "  +0.51%  "
Set-TextValue $ws.Range("D48") "86.24"
Set-TextValue $ws.Range("E48") "  +0.68%  "
Set-TextValue $ws.Range("D49") "0.0511"
Set-TextValue $ws.Range("E49") "  -0.10%  "
Set-TextValue $ws.Range("E50") "  +1.09%  "
Set-TextValue $ws.Range("B51") "BabyDogeCoin"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws.Range("D51") "0.0₇0960"
Set-TextValue $ws.Range("E51") "  -1.34%  "
